# Populate shipping address information / payment information
# Target sheet: "Address" (2nd worksheet) — rows 4/5 are the "Shipping Address"
# section mirroring the "Billing Address" section in rows 1/2, but the
# Area Code / Primary Phone fields are removed in favor of a single
# "Company Name" field, whose value moves from column H into column F.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Address")

# Make sure the Address sheet is the active one (it already is, but be explicit).
$ws.Activate()

# --- Shipping header row (row 4): drop "Area Code" / "Primary Phone" and
#     replace the "Area Code" header with "Company Name" (previously in H4).
$ws.Range("F4").Value = "Company Name"
$ws.Range("G4").Clear()
$ws.Range("H4").Clear()

# --- Shipping data row (row 5): drop the area code / phone number values
#     and move the company name ("Apple") from H5 into F5. H5 had no
#     explicit style applied, so clear F5's existing (numeric) style first.
$ws.Range("F5").Clear()
$ws.Range("F5").Style = "Normal"
$ws.Range("F5").Value = "Apple"
$ws.Range("G5").Clear()
$ws.Range("H5").Clear()

# --- Update the active selection on the sheet.
$ws.Range("G9").Select() | Out-Null
